$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 324987
$ws.Range("D2").Value = 414031423
$ws.Range("C3").Value = 262
$ws.Range("D3").Value = 313452
$ws.Range("C8").Value = 870
$ws.Range("D8").Value = 1279795
$ws.Range("C10").Value = 118155
$ws.Range("D10").Value = 173121141
$ws.Range("C12").Value = 60245
$ws.Range("D12").Value = 86947766
$ws.Range("C16").Value = 4031
$ws.Range("D16").Value = 5721197
$ws.Range("C20").Value = 6885
$ws.Range("D20").Value = 9610848
$ws.Range("C22").Value = 78348
$ws.Range("D22").Value = 97636176
$ws.Range("C28").Value = 32723
$ws.Range("D28").Value = 47895791
$ws.Range("C30").Value = 11600
$ws.Range("D30").Value = 16686097
$ws.Range("C33").Value = 1569
$ws.Range("D33").Value = 2204781
$ws.Range("C35").Value = 1892
$ws.Range("D35").Value = 2673168
$ws.Range("C36").Value = 98240
$ws.Range("D36").Value = 123565614
$ws.Range("C44").Value = 44695
$ws.Range("D44").Value = 65501066
$ws.Range("C46").Value = 9255
$ws.Range("D46").Value = 13276108
$ws.Range("C48").Value = 1418
$ws.Range("D48").Value = 1969311
$ws.Range("C51").Value = 2407
$ws.Range("D51").Value = 3363848
$ws.Range("C52").Value = 69835
$ws.Range("D52").Value = 87585741
$ws.Range("C59").Value = 28446
$ws.Range("D59").Value = 41716704
$ws.Range("C62").Value = 11262
$ws.Range("D62").Value = 16285835
$ws.Range("C68").Value = 1531
$ws.Range("D68").Value = 2144417
$ws.Range("C70").Value = 20715
$ws.Range("D70").Value = 27127758
$ws.Range("C74").Value = 7669
$ws.Range("D74").Value = 11230356
$ws.Range("C76").Value = 5186
$ws.Range("D76").Value = 7529986
$ws.Range("C79").Value = 142532
$ws.Range("D79").Value = 177627259
$ws.Range("C85").Value = 64128
$ws.Range("D85").Value = 93984195
$ws.Range("C86").Value = 82
$ws.Range("D86").Value = 121582
$ws.Range("C88").Value = 30089
$ws.Range("D88").Value = 43524018
$ws.Range("C90").Value = 2755
$ws.Range("D90").Value = 3966957
$ws.Range("C91").Value = 2917
$ws.Range("D91").Value = 4124547
$ws.Range("C92").Value = 33857
$ws.Range("D92").Value = 45886126
$ws.Range("C96").Value = 8212
$ws.Range("D96").Value = 12072127
$ws.Range("C98").Value = 7549
$ws.Range("D98").Value = 10951980
$ws.Range("C100").Value = 543
$ws.Range("D100").Value = 771135
$ws.Range("C102").Value = 10635
$ws.Range("D102").Value = 16360661
$ws.Range("C104").Value = 2610
$ws.Range("D104").Value = 4285676
$ws.Range("C106").Value = 3539
$ws.Range("D106").Value = 5817151
$ws.Range("C109").Value = 208
$ws.Range("D109").Value = 326860
$ws.Range("C110").Value = 143179
$ws.Range("D110").Value = 177042076
$ws.Range("C114").Value = 954
$ws.Range("D114").Value = 1399477
$ws.Range("C116").Value = 53212
$ws.Range("D116").Value = 77990734
$ws.Range("C118").Value = 27495
$ws.Range("D118").Value = 39835838
$ws.Range("C122").Value = 2328
$ws.Range("D122").Value = 3270699
$ws.Range("C124").Value = 523430
$ws.Range("D124").Value = 691300140
$ws.Range("C126").Value = 219
$ws.Range("D126").Value = 322509
$ws.Range("C129").Value = 1391
$ws.Range("D129").Value = 2061682
$ws.Range("C131").Value = 211149
$ws.Range("D131").Value = 310392931
$ws.Range("C132").Value = 411
$ws.Range("D132").Value = 613250
$ws.Range("C134").Value = 186797
$ws.Range("D134").Value = 271646081
$ws.Range("C137").Value = 2873
$ws.Range("D137").Value = 4034299
$ws.Range("C139").Value = 6567
$ws.Range("D139").Value = 9278534
$ws.Range("C142").Value = 45408
$ws.Range("D142").Value = 60617729
$ws.Range("C148").Value = 14244
$ws.Range("D148").Value = 20883723
$ws.Range("C149").Value = 3815
$ws.Range("D149").Value = 5501740
$ws.Range("C154").Value = 398
$ws.Range("D154").Value = 562763
$ws.Range("C155").Value = 17900
$ws.Range("D155").Value = 23657680
$ws.Range("C159").Value = 7289
$ws.Range("D159").Value = 10605404
$ws.Range("C161").Value = 5095
$ws.Range("D161").Value = 7334996
$ws.Range("C163").Value = 282
$ws.Range("D163").Value = 390239
$ws.Range("C166").Value = 18907
$ws.Range("D166").Value = 30982945
$ws.Range("C167").Value = 2054
$ws.Range("D167").Value = 3385072
$ws.Range("C171").Value = 107
$ws.Range("D171").Value = 182949
$ws.Range("C172").Value = 88804
$ws.Range("D172").Value = 110979784
$ws.Range("C177").Value = 646
$ws.Range("D177").Value = 952088
$ws.Range("C179").Value = 34175
$ws.Range("D179").Value = 50112703
$ws.Range("C181").Value = 13193
$ws.Range("D181").Value = 19061357
$ws.Range("C183").Value = 1253
$ws.Range("D183").Value = 1754429
$ws.Range("C185").Value = 1707
$ws.Range("D185").Value = 2397239
$ws.Range("C187").Value = 241175
$ws.Range("D187").Value = 299679902
$ws.Range("C189").Value = 171
$ws.Range("D189").Value = 246736
$ws.Range("C193").Value = 884
$ws.Range("D193").Value = 1300345
$ws.Range("C195").Value = 87318
$ws.Range("D195").Value = 127993230
$ws.Range("C198").Value = 33420
$ws.Range("D198").Value = 48109725
$ws.Range("C201").Value = 5161
$ws.Range("D201").Value = 7347296
$ws.Range("C204").Value = 5038
$ws.Range("D204").Value = 6976710
$ws.Range("C207").Value = 267233
$ws.Range("D207").Value = 330666003
$ws.Range("C209").Value = 258
$ws.Range("D209").Value = 368539
$ws.Range("C216").Value = 96065
$ws.Range("D216").Value = 140534619
$ws.Range("C219").Value = 52144
$ws.Range("D219").Value = 75364853
$ws.Range("C222").Value = 4716
$ws.Range("D222").Value = 6620323
$ws.Range("C225").Value = 5985
$ws.Range("D225").Value = 8283726
$ws.Range("C228").Value = 107639
$ws.Range("D228").Value = 134553886
$ws.Range("C230").Value = 75
$ws.Range("D230").Value = 107513
$ws.Range("C235").Value = 49954
$ws.Range("D235").Value = 73179763
$ws.Range("C237").Value = 12626
$ws.Range("D237").Value = 18154490
$ws.Range("C241").Value = 2592
$ws.Range("D241").Value = 3628686
$ws.Range("C242").Value = 261211
$ws.Range("D242").Value = 329794298
$ws.Range("C243").Value = 173
$ws.Range("D243").Value = 214433
$ws.Range("C250").Value = 96780
$ws.Range("D250").Value = 141808870
$ws.Range("C253").Value = 66067
$ws.Range("D253").Value = 95744711
$ws.Range("C255").Value = 2435
$ws.Range("D255").Value = 3433773
$ws.Range("C258").Value = 4760
$ws.Range("D258").Value = 6690137
